$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of C4: "More flexible estimation methods/ management procedures"
# becomes "More flexible management procedures".
$ws.Range("C4").Value = "More flexible management procedures"

# Make C4 the active/selected cell (as reflected by the new <selection> in sheetView).
$ws.Range("C4").Select()
